$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Small word-fixes inside the existing "- Quem sou eu, ..." paragraph
# ------------------------------------------------------------------
# "por que" -> "porquê"
$d.Content.Find.Execute("por que", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "porquê", 2) | Out-Null

# "nos enviado." -> "nos enviados."
$d.Content.Find.Execute("nos enviado.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "nos enviados.", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Grow the document with a blank line followed by two new narrative
#    paragraphs about the "Criação de Jogos" theme.
# ------------------------------------------------------------------
$d.Paragraphs.Last.Range.InsertParagraphAfter()   # blank line
$d.Paragraphs.Last.Range.InsertParagraphAfter()   # placeholder for 1st new paragraph
$d.Paragraphs.Last.Range.InsertParagraphAfter()   # placeholder for 2nd new paragraph

$d.Paragraphs(5).Range.Text = "Meu tema é Criação de Jogos e eu me chamo Eduardo Ribeiro Santos Nascimento, tenho 19 anos e atualmente estou na SPTECH, mas como eu cheguei aqui e o que este tema tem a ver comigo?"

$d.Paragraphs(6).Range.Text = "Ele se relaciona comigo desde quando eu era pequeno e nem sabia o que era criar um jogo. Eu gostava de brincar de carrinho e boneco, dentro dessas brincadeiras eu imaginava e desenvolvia histórias"
